$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Name = "၂၀၀၄-၂၀၁၃"
